$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The worksheet is protected; unprotect it so the cells below can be edited.
$ws.Unprotect()

# Update the "as of" date embedded in the confidentiality disclaimer text (A7)
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."
$ws.Rows("7").AutoFit()

# Update Weight (D) and Percent Change (E) values for rows 2-4
$ws.Range("D2").Value = 0.8474552175149276
$ws.Range("E2").Value = 0.002264720684448829

$ws.Range("D3").Value = 0.1525447824850725
$ws.Range("E3").Value = 0.004287045666356182

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.002573215808927998

# Restore sheet protection (best-effort; original used a legacy password hash
# that cannot be re-derived, so we reapply protection without a password).
$ws.Protect()
